$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7655548453330994
$ws.Range("B1").Value = 1.332211017608643
$ws.Range("C1").Value = 4.26425838470459
$ws.Range("D1").Value = 3.827998399734497
$ws.Range("E1").Value = 1.673259735107422
